# Weekly update: two new "Ajo" price records were captured for
# Terminal Hortofrutícola Agro Chillán (week of 2023-03-28, serial 45013).
# They are inserted right after the existing row for 2021-11-30 (row 350),
# which pushes every subsequent record down by two rows; the two oldest
# records that fall off the bottom of the block are preserved by being
# re-appended as the new last two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows, shifting rows 351-376 down to 353-378
# (Excel's own Insert() carries formatting/shared content along, matching
# what the diff shows for the shifted rows).
$ws.Range("A351:R352").Insert()

# New row 351: $/caja 10 kilos record
$ws.Range("A351").Value = 7
$ws.Range("B351").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C351").Value = "Ñuble"
$ws.Range("D351").Value = 45013
$ws.Range("D351").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E351").Value = 16
$ws.Range("F351").Value = 100112003
$ws.Range("G351").Value = "Ajo"
$ws.Range("H351").Value = "Chino"
$ws.Range("I351").Value = "Primera"
$ws.Range("J351").Value = 80
$ws.Range("K351").Value = 17000
$ws.Range("L351").Value = 18000
$ws.Range("M351").Value = 17500
$ws.Range("N351").Value = "$/caja 10 kilos"
$ws.Range("O351").Value = "China"
$ws.Range("P351").Value = 1750
$ws.Range("Q351").Value = 10
$ws.Range("R351").Value = "Hortaliza"

# New row 352: $/malla 10 kilos record, same date
$ws.Range("A352").Value = 7
$ws.Range("B352").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C352").Value = "Ñuble"
$ws.Range("D352").Value = 45013
$ws.Range("D352").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E352").Value = 16
$ws.Range("F352").Value = 100112003
$ws.Range("G352").Value = "Ajo"
$ws.Range("H352").Value = "Chino"
$ws.Range("I352").Value = "Primera"
$ws.Range("J352").Value = 60
$ws.Range("K352").Value = 20000
$ws.Range("L352").Value = 20000
$ws.Range("M352").Value = 20000
$ws.Range("N352").Value = "$/malla 10 kilos"
$ws.Range("O352").Value = "China"
$ws.Range("P352").Value = 2000
$ws.Range("Q352").Value = 10
$ws.Range("R352").Value = "Hortaliza"
